$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 0.8944602796882041
$ws.Range("C2").Value = 0.8944602796882039
$ws.Range("D2").Value = 0.8944602796882039

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.9991192367107327
$ws.Range("C3").Value = 0.9987285627517402
$ws.Range("D3").Value = 0.9806305037591442

# Row 4 - name changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9987837031296761
$ws.Range("C4").Value = 0.9981038032102866
$ws.Range("D4").Value = 0.994749339258567

# Row 5 - name changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9994027109757438
$ws.Range("C5").Value = 0.9993431277534596
$ws.Range("D5").Value = 0.9990875095672079
